$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.934691846370697
$ws.Range("B1").Value = 0.9854012131690979
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 2.17516040802002
$ws.Range("E1").Value = 1.075182318687439
